$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 172, pushing existing rows 172:185 down to 173:186.
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row 172 with its data.
$ws.Range("A172").Value = 11
$ws.Range("B172").Value = "Vega Monumental Concepción"
$ws.Range("C172").Value = "Bíobío"
$ws.Range("D172").Value = 45147
$ws.Range("E172").Value = 8
$ws.Range("F172").Value = "Fruta"
$ws.Range("G172").Value = 100108
$ws.Range("H172").Value = "Tropicales y subtropicales"
$ws.Range("I172").Value = 100108002
$ws.Range("J172").Value = "Mango"
$ws.Range("K172").Value = "Sin especificar"
$ws.Range("L172").Value = "Primera"
$ws.Range("M172").Value = 180
$ws.Range("N172").Value = 7000
$ws.Range("O172").Value = 7500
$ws.Range("P172").Value = 7278
$ws.Range("Q172").Value = "$/bandeja 4 kilos"
$ws.Range("R172").Value = "Brasil"
$ws.Range("S172").Value = 1820
$ws.Range("T172").Value = 4
